$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.899.05'
$ws.Range('E2').Value = '  -1.30%  '

$ws.Range('D3').Value = '2.541.93'
$ws.Range('E3').Value = '  -0.03%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.08%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.79'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.46%  '

$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.582'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.22%  '

$ws.Range('E9').Value = '  -1.52%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.48'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.80%  '

$ws.Range('E11').Value = '  -0.56%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.354'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.19%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.26'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.05%  '

$ws.Range('D14').Value = '2.997.33'
$ws.Range('E14').Value = '  -0.23%  '

$ws.Range('D15').Value = '62.865.51'
$ws.Range('E15').Value = '  -1.10%  '

$ws.Range('E16').Value = '  -0.99%  '

$ws.Range('D17').Value = '2.564.47'
$ws.Range('E17').Value = '  +0.61%  '

$ws.Range('E18').Value = '  -1.97%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '335.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.89%  '

$ws.Range('E20').Value = '  -0.81%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.16%  '

$ws.Range('E22').Value = '  +0.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.54%  '

$ws.Range('E24').Value = '  -0.36%  '

$ws.Range('E25').Value = '  +1.07%  '

$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('E27').Value = '  -0.25%  '

$ws.Range('E28').Value = '  +1.68%  '

$ws.Range('E29').Value = '  +4.80%  '

$ws.Range('D30').Value = '0.0₃0809'
$ws.Range('E30').Value = '  -2.67%  '

$ws.Range('E31').Value = '  -1.21%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '177.67'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.41%  '

$ws.Range('E33').Value = '  -3.74%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '402.61'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.77%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '19.08'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.19%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.398'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.41%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.21%  '

$ws.Range('E39').Value = '  -1.75%  '

$ws.Range('E40').Value = '  +0.01%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '39.26'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.33%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '150.80'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.35%  '

$ws.Range('E43').Value = '  -1.74%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.74'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.09%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0534'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.42%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.598'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.45%  '

$ws.Range('E47').Value = '  -0.68%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0237'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.06%  '

$ws.Range('E49').Value = '  -3.56%  '

$ws.Range('E51').Value = '  -7.58%  '
